# Implement Runmode for Test suite level and test case level
#
# 1. Add a new "test_suite" sheet (run-level Runmode control) as the first
#    tab, ahead of the existing AddCustomerTest / OpenAccountTest sheets.
# 2. Add a test-case level "runmode" column to AddCustomerTest, and rename
#    its old "alerttest" header to "alerttext".

$wb = $excel.ActiveWorkbook

# --- 1. New "test_suite" sheet, inserted before the first existing sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$suite = $wb.Worksheets.Add($firstSheet)
$suite.Name = "test_suite"

$suite.Range("A1").Value = "TCID"
$suite.Range("B1").Value = "Runmode"
$suite.Range("A2").Value = "BankManagerLoginTest"
$suite.Range("B2").Value = "Y"
$suite.Range("A3").Value = "AddCustomerTest"
$suite.Range("B3").Value = "Y"
$suite.Range("A4").Value = "OpenAccountTest"
$suite.Range("B4").Value = "Y"
$suite.Range("B4").Select()

# --- 2. Test-case level Runmode column on AddCustomerTest ---
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")

$addCustomer.Range("E3").Value = "N"
$addCustomer.Range("E1").Value = "runmode"
$addCustomer.Range("E2").Value = "Y"
$addCustomer.Range("E4").Value = "Y"
$addCustomer.Range("D1").Value = "alerttext"
$addCustomer.Range("D1").Select()
